$p = $ppt.ActivePresentation

# --- Slide 1: update the "Generated on" date subtitle ---
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(2)
$subtitle.TextFrame.TextRange.Text = "Generated on 2026-02-13"

# --- Slide 2: update default paragraph spacing on the bullet list ---
$slide2 = $p.Slides.Item(2)
$body = $slide2.Shapes.Item(3)
$tr = $body.TextFrame.TextRange
$count = $tr.Paragraphs().Count

for ($i = 1; $i -le $count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.IndentLevel = 2
    $para.ParagraphFormat.SpaceBefore = 3
    $para.ParagraphFormat.SpaceAfter = 3
}
